# Updated blog_template.docx with corrected placeholders
#
# The template's first paragraph currently holds a merge-field-style run
# for {{TITLE}} (fldChar begin/instrText/separate/end). The corrected
# template needs two new leading paragraphs that hold the placeholders as
# plain literal text:
#   1) a new paragraph containing the literal text "{{TITLE}}"
#   2) a new paragraph containing the literal text "{{CONTENT}}" followed
#      by (i.e. sharing a paragraph with) whatever used to be the very
#      first paragraph's content (left completely untouched).
#
# We do this with a single surgical Range.InsertXML call collapsed to the
# very start of the document: inserting two fresh <w:p> elements there
# pushes the original first paragraph's existing runs (the fldChar field
# code for {{TITLE}}) down into the second of the two new paragraphs,
# without disturbing those original runs at all.

$d = $word.ActiveDocument

$firstParagraph = $d.Paragraphs.Item(1)
$insertionPoint = $firstParagraph.Range.Duplicate
$insertionPoint.Collapse(1)   # wdCollapseStart

$xmlFragment = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p><w:r><w:t>{{TITLE}}</w:t></w:r></w:p>' +
            '<w:p><w:r><w:t>{{CONTENT}}</w:t></w:r></w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$insertionPoint.InsertXML($xmlFragment)
